{"js": "// Update the multiplication-problem answers throughout the table.\n// Each cell holds a single run of text like \"659\u00d77=4613\"; we locate the\n// old expression via a body search and replace it with the new one.\nconst replacements = [\n  [\"659\u00d77=4613\", \"456\u00d76=2736\"],\n  [\"292\u00d74=1168\", \"651\u00d79=5859\"],\n  [\"616\u00d79=5544\", \"540\u00d77=3780\"],\n  [\"296\u00d78=2368\", \"646\u00d78=5168\"],\n  [\"850\u00d74=3400\", \"132\u00d74=528\"],\n  [\"946\u00d77=6622\", \"905\u00d74=3620\"],\n  [\"566\u00d75=2830\", \"979\u00d73=2937\"],\n  [\"818\u00d77=5726\", \"197\u00d78=1576\"],\n  [\"486\u00d72=972\", \"193\u00d72=386\"],\n  [\"502\u00d78=4016\", \"225\u00d78=1800\"],\n  [\"542\u00d72=1084\", \"902\u00d72=1804\"],\n  [\"389\u00d74=1556\", \"685\u00d78=5480\"],\n  [\"502\u00d73=1506\", \"231\u00d78=1848\"],\n  [\"446\u00d76=2676\", \"180\u00d74=720\"],\n  [\"515\u00d75=2575\", \"739\u00d72=1478\"],\n  [\"774\u00d78=6192\", \"326\u00d72=652\"],\n  [\"774\u00d73=2322\", \"255\u00d73=765\"],\n  [\"862\u00d78=6896\", \"274\u00d78=2192\"],\n  [\"981\u00d78=7848\", \"709\u00d77=4963\"],\n  [\"129\u00d79=1161\", \"867\u00d76=5202\"],\n  [\"724\u00d77=5068\", \"848\u00d74=3392\"],\n  [\"964\u00d77=6748\", \"114\u00d77=798\"],\n  [\"459\u00d79=4131\", \"425\u00d75=2125\"],\n  [\"656\u00d77=4592\", \"586\u00d74=2344\"],\n  [\"718\u00d72=1436\", \"864\u00d72=1728\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication-problem answers throughout the table.\n# Each cell holds a single run of text like \"659\u00d77=4613\"; we use\n# Find/Replace (Replace All) against the whole document content for\n# each old/new pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"659\u00d77=4613\", \"456\u00d76=2736\"),\n    @(\"292\u00d74=1168\", \"651\u00d79=5859\"),\n    @(\"616\u00d79=5544\", \"540\u00d77=3780\"),\n    @(\"296\u00d78=2368\", \"646\u00d78=5168\"),\n    @(\"850\u00d74=3400\", \"132\u00d74=528\"),\n    @(\"946\u00d77=6622\", \"905\u00d74=3620\"),\n    @(\"566\u00d75=2830\", \"979\u00d73=2937\"),\n    @(\"818\u00d77=5726\", \"197\u00d78=1576\"),\n    @(\"486\u00d72=972\", \"193\u00d72=386\"),\n    @(\"502\u00d78=4016\", \"225\u00d78=1800\"),\n    @(\"542\u00d72=1084\", \"902\u00d72=1804\"),\n    @(\"389\u00d74=1556\", \"685\u00d78=5480\"),\n    @(\"502\u00d73=1506\", \"231\u00d78=1848\"),\n    @(\"446\u00d76=2676\", \"180\u00d74=720\"),\n    @(\"515\u00d75=2575\", \"739\u00d72=1478\"),\n    @(\"774\u00d78=6192\", \"326\u00d72=652\"),\n    @(\"774\u00d73=2322\", \"255\u00d73=765\"),\n    @(\"862\u00d78=6896\", \"274\u00d78=2192\"),\n    @(\"981\u00d78=7848\", \"709\u00d77=4963\"),\n    @(\"129\u00d79=1161\", \"867\u00d76=5202\"),\n    @(\"724\u00d77=5068\", \"848\u00d74=3392\"),\n    @(\"964\u00d77=6748\", \"114\u00d77=798\"),\n    @(\"459\u00d79=4131\", \"425\u00d75=2125\"),\n    @(\"656\u00d77=4592\", \"586\u00d74=2344\"),\n    @(\"718\u00d72=1436\", \"864\u00d72=1728\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $oldText\n    $range.Find.Replacement.Text = $newText\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
